$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text parses as a plain number (e.g. "11.60") need an
# explicit text format first, otherwise Excel stores them as numeric values and
# drops the trailing zero / trailing precision (matches the source data, which is
# plain text scraped from the site, e.g. "227.75", "11.60", "6.00").
$textCells = @("D5", "D6", "D8", "D9", "D10", "D13", "D15", "D16", "D18", "D19", "D21", "D22", "D27", "D28", "D32", "D34", "D36", "D37", "D39", "D44", "D46", "D49", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "35.384.45"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "1.846.86"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "227.75"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "42.14"
$ws.Range("E8").Value = "  +11.43%  "
$ws.Range("D9").Value = "0.308"
$ws.Range("E9").Value = "  +5.31%  "
$ws.Range("D10").Value = "0.0688"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("D12").Value = "2.113.17"
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("D13").Value = "11.60"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").Value = "1.842.02"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "4.74"
$ws.Range("E15").Value = "  +6.54%  "
$ws.Range("D16").Value = "0.665"
$ws.Range("E16").Value = "  +4.55%  "
$ws.Range("D17").Value = "35.304.56"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "70.14"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "245.36"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "0.0₃0793"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "12.18"
$ws.Range("E21").Value = "  +8.14%  "
$ws.Range("D22").Value = "4.80"
$ws.Range("E22").Value = "  +15.77%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "17.88"
$ws.Range("E27").Value = "  +2.84%  "
$ws.Range("D28").Value = "0.123"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "3.570.96"
$ws.Range("E29").Value = "  +46.97%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +7.11%  "
$ws.Range("D32").Value = "3.94"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").Value = "0.0534"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D36").Value = "0.676"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").Value = "88.91"
$ws.Range("E37").Value = "  +9.39%  "
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").Value = "1.03"
$ws.Range("E39").Value = "  +9.29%  "
$ws.Range("D40").Value = "1.336.47"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("D44").Value = "14.92"
$ws.Range("E44").Value = "  +5.15%  "
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "2.83"
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("D48").Value = "2.012.21"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "6.00"
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("D50").Value = "104.09"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  +0.02%  "
